$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-08"

# Update the header label cell (I1) that shows "2022 (through 11-07)"
$ws.Range("I1").Value = "2022 (through 11-08)"

# Update the data values that changed for the new date
$ws.Range("I9").Value = 162
$ws.Range("I12").Value = 22
$ws.Range("I14").Value = 1421
